$d = $word.ActiveDocument

# Replace the entire first paragraph's text (which spans multiple runs
# separated by manual line breaks) with the single word "jjjjj".
$para = $d.Paragraphs(1)
$range = $para.Range
# Trim the trailing paragraph mark from the range so we only replace the
# visible text/content, not the paragraph mark itself.
$range.MoveEnd(1, -1) | Out-Null
$range.Text = "jjjjj"

Write-Output "done"
